$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H - Labor Booking User
$ws.Range("H1").Value = "Labor Booking User"
$ws.Range("H2").Value = "a811K0000004fpN"
$ws.Range("H3").Value = "a811K0000004fpN"

# Header row remainder - I1:K1
$ws.Range("I1").Value = "SiteID"
$ws.Range("J1").Value = "Location ID"
$ws.Range("K1").Value = "Location Number"

# Rows 2-3, columns I:K
$ws.Range("I2").Value = "a7q410000004I1W"
$ws.Range("J2").Value = "a7Z4100000000hb"
$ws.Range("K2").Value = "SY_ReceiptLoc"
$ws.Range("I3").Value = "a7q410000004I1W"
$ws.Range("J3").Value = "a7Z4100000000hb"
$ws.Range("K3").Value = "SY_ReceiptLoc"

# Column width changes (values chosen so the resulting stored width matches
# the target as closely as this engine's ColumnWidth->stored-width rounding allows)
$ws.Range("B1").ColumnWidth = 20.59
$ws.Range("H1").ColumnWidth = 16.09
$ws.Range("I1").ColumnWidth = 15.59
$ws.Range("J1").ColumnWidth = 15.09
$ws.Range("K1").ColumnWidth = 14.25

# Selection change
[void]$ws.Range("H1:H1048576").Select()
